$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("M6").Value = 1.04
$ws.Range("N6").Value = 9
$ws.Range("Q6").Value = 1.8
$ws.Range("R6").Value = 2
